$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.423576
$ws.Range("H2").Value = 19.270728
$ws.Range("I2").Value = 0.001681024218962088
$ws.Range("J2").Value = 0.001681024218962088
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.672731
$ws.Range("N2").Value = 8.018193
$ws.Range("O2").Value = 0.5408378022089502
$ws.Range("P2").Value = 0.5408378022089503
$ws.Range("Q2").Value = 17.168490706056
$ws.Range("R2").Value = 154.516416354504
$ws.Range("S2").Value = 0.0009091614440434728
$ws.Range("T2").Value = 0.000909161444043473

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.423576
$ws.Range("H3").Value = 19.270728
$ws.Range("I3").Value = 0.001681024218962088
$ws.Range("J3").Value = 0.001681024218962088
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.627877666666667
$ws.Range("N3").Value = 4.883633
$ws.Range("O3").Value = 0.3294075533620981
$ws.Range("P3").Value = 0.3294075533620982
$ws.Range("Q3").Value = 10.456795910536
$ws.Range("R3").Value = 94.11116319482399
$ws.Range("S3").Value = 0.0005537420751107334
$ws.Range("T3").Value = 0.0005537420751107335

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.423576
$ws.Range("H4").Value = 19.270728
$ws.Range("I4").Value = 0.001681024218962088
$ws.Range("J4").Value = 0.001681024218962088
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01852966666666667
$ws.Range("N4").Value = 0.055589
$ws.Range("O4").Value = 0.003749552123152102
$ws.Range("P4").Value = 0.003749552123152104
$ws.Range("Q4").Value = 0.119026722088
$ws.Range("R4").Value = 1.071240498792
$ws.Range("S4").Value = 0.000006303087929279403
$ws.Range("T4").Value = 0.000006303087929279404

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.423576
$ws.Range("H5").Value = 19.270728
$ws.Range("I5").Value = 0.001681024218962088
$ws.Range("J5").Value = 0.001681024218962088
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6226963333333334
$ws.Range("N5").Value = 1.868089
$ws.Range("O5").Value = 0.1260050923057995
$ws.Range("P5").Value = 0.1260050923057995
$ws.Range("Q5").Value = 3.999937222088001
$ws.Range("R5").Value = 35.999434998792
$ws.Range("S5").Value = 0.0002118176118786024
$ws.Range("T5").Value = 0.0002118176118786025

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3580.644531333333
$ws.Range("H6").Value = 10741.933594
$ws.Range("I6").Value = 0.9370403925578976
$ws.Range("J6").Value = 0.9370403925578976
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.672731
$ws.Range("N6").Value = 8.018193
$ws.Range("O6").Value = 0.5408378022089502
$ws.Range("P6").Value = 0.5408378022089503
$ws.Range("Q6").Value = 9570.099638875072
$ws.Range("R6").Value = 86130.89674987564
$ws.Range("S6").Value = 0.5067868664920252
$ws.Range("T6").Value = 0.5067868664920253

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3580.644531333333
$ws.Range("H7").Value = 10741.933594
$ws.Range("I7").Value = 0.9370403925578976
$ws.Range("J7").Value = 0.9370403925578976
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.627877666666667
$ws.Range("N7").Value = 4.883633
$ws.Range("O7").Value = 0.3294075533620981
$ws.Range("P7").Value = 0.3294075533620982
$ws.Range("Q7").Value = 5828.851264829666
$ws.Range("R7").Value = 52459.661383467
$ws.Range("S7").Value = 0.308668183113957
$ws.Range("T7").Value = 0.3086681831139571

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3580.644531333333
$ws.Range("H8").Value = 10741.933594
$ws.Range("I8").Value = 0.9370403925578976
$ws.Range("J8").Value = 0.9370403925578976
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01852966666666667
$ws.Range("N8").Value = 0.055589
$ws.Range("O8").Value = 0.003749552123152102
$ws.Range("P8").Value = 0.003749552123152104
$ws.Range("Q8").Value = 66.34814961742956
$ws.Range("R8").Value = 597.1333465568659
$ws.Range("S8").Value = 0.003513481793394745
$ws.Range("T8").Value = 0.003513481793394746

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3580.644531333333
$ws.Range("H9").Value = 10741.933594
$ws.Range("I9").Value = 0.9370403925578976
$ws.Range("J9").Value = 0.9370403925578976
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6226963333333334
$ws.Range("N9").Value = 1.868089
$ws.Range("O9").Value = 0.1260050923057995
$ws.Range("P9").Value = 0.1260050923057995
$ws.Range("Q9").Value = 2229.654220631319
$ws.Range("R9").Value = 20066.88798568187
$ws.Range("S9").Value = 0.1180718611585205
$ws.Range("T9").Value = 0.1180718611585205

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col1a2"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9157713333333334
$ws.Range("H10").Value = 2.747314
$ws.Range("I10").Value = 0.0002396537054071653
$ws.Range("J10").Value = 0.0002396537054071653
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 2.672731
$ws.Range("N10").Value = 8.018193
$ws.Range("O10").Value = 0.5408378022089502
$ws.Range("P10").Value = 0.5408378022089503
$ws.Range("Q10").Value = 2.447610431511334
$ws.Range("R10").Value = 22.028493883602
$ws.Range("S10").Value = 0.0001296137833236425
$ws.Range("T10").Value = 0.0001296137833236425

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col1a2"
$ws.Range("C11").Value = "Itga2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.9157713333333334
$ws.Range("H11").Value = 2.747314
$ws.Range("I11").Value = 0.0002396537054071653
$ws.Range("J11").Value = 0.0002396537054071653
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.627877666666667
$ws.Range("N11").Value = 4.883633
$ws.Range("O11").Value = 0.3294075533620981
$ws.Range("P11").Value = 0.3294075533620982
$ws.Range("Q11").Value = 1.490763701306889
$ws.Range("R11").Value = 13.416873311762
$ws.Range("S11").Value = 0.00007894374075233533
$ws.Range("T11").Value = 0.00007894374075233536

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col1a2"
$ws.Range("C12").Value = "Itga2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9157713333333334
$ws.Range("H12").Value = 2.747314
$ws.Range("I12").Value = 0.0002396537054071653
$ws.Range("J12").Value = 0.0002396537054071653
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01852966666666667
$ws.Range("N12").Value = 0.055589
$ws.Range("O12").Value = 0.003749552123152102
$ws.Range("P12").Value = 0.003749552123152104
$ws.Range("Q12").Value = 0.01696893754955556
$ws.Range("R12").Value = 0.152720437946
$ws.Range("S12").Value = 0.000000898594059930705
$ws.Range("T12").Value = 0.0000008985940599307054

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col1a2"
$ws.Range("C13").Value = "Itga2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9157713333333334
$ws.Range("H13").Value = 2.747314
$ws.Range("I13").Value = 0.0002396537054071653
$ws.Range("J13").Value = 0.0002396537054071653
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6226963333333334
$ws.Range("N13").Value = 1.868089
$ws.Range("O13").Value = 0.1260050923057995
$ws.Range("P13").Value = 0.1260050923057995
$ws.Range("Q13").Value = 0.5702474514384446
$ws.Range("R13").Value = 5.132227062946001
$ws.Range("S13").Value = 0.00003019758727125674
$ws.Range("T13").Value = 0.00003019758727125675

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col1a2"
$ws.Range("C14").Value = "Itga2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 233.243637
$ws.Range("H14").Value = 699.7309110000001
$ws.Range("I14").Value = 0.0610389295177331
$ws.Range("J14").Value = 0.06103892951773311
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 2.672731
$ws.Range("N14").Value = 8.018193
$ws.Range("O14").Value = 0.5408378022089502
$ws.Range("P14").Value = 0.5408378022089503
$ws.Range("Q14").Value = 623.3974991626471
$ws.Range("R14").Value = 5610.577492463824
$ws.Range("S14").Value = 0.03301216048955778
$ws.Range("T14").Value = 0.0330121604895578

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col1a2"
$ws.Range("C15").Value = "Itga2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 233.243637
$ws.Range("H15").Value = 699.7309110000001
$ws.Range("I15").Value = 0.0610389295177331
$ws.Range("J15").Value = 0.06103892951773311
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.627877666666667
$ws.Range("N15").Value = 4.883633
$ws.Range("O15").Value = 0.3294075533620981
$ws.Range("P15").Value = 0.3294075533620982
$ws.Range("Q15").Value = 379.692107564407
$ws.Range("R15").Value = 3417.228968079663
$ws.Range("S15").Value = 0.02010668443227801
$ws.Range("T15").Value = 0.02010668443227802

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col1a2"
$ws.Range("C16").Value = "Itga2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 233.243637
$ws.Range("H16").Value = 699.7309110000001
$ws.Range("I16").Value = 0.0610389295177331
$ws.Range("J16").Value = 0.06103892951773311
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.01852966666666667
$ws.Range("N16").Value = 0.055589
$ws.Range("O16").Value = 0.003749552123152102
$ws.Range("P16").Value = 0.003749552123152104
$ws.Range("Q16").Value = 4.321926845731
$ws.Range("R16").Value = 38.89734161157901
$ws.Range("S16").Value = 0.0002288686477681477
$ws.Range("T16").Value = 0.0002288686477681478

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col1a2"
$ws.Range("C17").Value = "Itga2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 233.243637
$ws.Range("H17").Value = 699.7309110000001
$ws.Range("I17").Value = 0.0610389295177331
$ws.Range("J17").Value = 0.06103892951773311
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6226963333333334
$ws.Range("N17").Value = 1.868089
$ws.Range("O17").Value = 0.1260050923057995
$ws.Range("P17").Value = 0.1260050923057995
$ws.Range("Q17").Value = 145.239957533231
$ws.Range("R17").Value = 1307.159617799079
$ws.Range("S17").Value = 0.007691215948129149
$ws.Range("T17").Value = 0.007691215948129152

